$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tenoxicam")

# Update the PharmGKB ID column value for row 2 (Tenoxicam / RxNorm / RxCUI row)
# from the text "C0076096" to the numeric value 37790
$ws.Range("D2").Value = 37790

# Update the selected cell on the sheet
$ws.Range("C10").Select()
